# Fix chatbot RAG flow and deck versioning
# Update the "Compute Metrics and Draft Summary" bullet list on slide 2
# (Content Placeholder 2) to reflect corrected WALT / rent calculations.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# --- Total GLA (paragraph 3): drop the per-tenant breakdown -----------------
$tr.Paragraphs(3, 1).Text = ""
$tr.Paragraphs(3, 1).Text = "**Total GLA**: 313,219 m" + [char]0x00B2

# --- Occupancy (paragraph 4): simplify wording ------------------------------
$tr.Paragraphs(4, 1).Text = ""
$tr.Paragraphs(4, 1).Text = "**Occupancy**: Fully leased"

# --- WALT (paragraph 5): rewrite + append the per-tenant breakdown ----------
$tr.Paragraphs(5, 1).Text = ""
$tr.Paragraphs(5, 1).Text = "**WALT (Weighted Average Lease Term)**: Calculation needed based on lease terms:"
[void]$tr.Paragraphs(5, 1).InsertAfter("`rIngram Micro: Lease ended in 2020 (0 years remaining)`rCNH Industrial: Lease ended in 2021 (0 years remaining)`rWALT = 0 years (as of the current date, assuming no renewals)")

# --- In-Place Rent (now paragraph 9): rewrite + append per-tenant rents ----
$tr.Paragraphs(9, 1).Text = ""
$tr.Paragraphs(9, 1).Text = "**In-Place Rent**:"
[void]$tr.Paragraphs(9, 1).InsertAfter("`rIngram Micro: " + [char]0x00A3 + "5.5/m" + [char]0x00B2 + "/year`rCNH Industrial: " + [char]0x00A3 + "4.5/m" + [char]0x00B2 + "/year")

# --- Key Highlight 1 (now paragraph 12) -------------------------------------
$tr.Paragraphs(12, 1).Text = ""
$tr.Paragraphs(12, 1).Text = "**Key Highlight 1**: Stable tenant base with major tenants like Ingram Micro and CNH Industrial."

# --- Key Highlight 2 (now paragraph 13) -------------------------------------
$tr.Paragraphs(13, 1).Text = ""
$tr.Paragraphs(13, 1).Text = "**Key Highlight 2**: Potential vacancy risk due to expired leases, requiring attention for renewals or new tenants."
